$d = $word.ActiveDocument

# ---------------------------------------------------------------
# Edit 1: "NAME: KAUSHIK NARAYANAN V" -> "NAME: " + "Ajay Kumar J"
#         (split into two runs with identical run formatting)
# ---------------------------------------------------------------
$d.Content.Find.Execute("NAME: KAUSHIK NARAYANAN V", $true, $false, $false, $false, $false, `
    $true, 1, $false, "NAME: ", 2)

$r1 = $d.Content
$r1.Find.Execute("NAME: ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r1.Collapse(0)
$r1.InsertAfter("Ajay Kumar J")
$r1.Font.Name = "Times New Roman"
$r1.Font.NameAscii = "Times New Roman"
$r1.Font.NameBi = "Times New Roman"
$r1.Font.Bold = $true
$r1.Font.Size = 14

# ---------------------------------------------------------------
# Edit 2: "REG NO: 192321047" -> "REG NO: 1923" + "72052"
#         (split into two runs with identical run formatting)
# ---------------------------------------------------------------
$d.Content.Find.Execute("REG NO: 192321047", $true, $false, $false, $false, $false, `
    $true, 1, $false, "REG NO: 1923", 2)

$r2 = $d.Content
$r2.Find.Execute("REG NO: 1923", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r2.Collapse(0)
$r2.InsertAfter("72052")
$r2.Font.Name = "Times New Roman"
$r2.Font.NameAscii = "Times New Roman"
$r2.Font.NameBi = "Times New Roman"
$r2.Font.Bold = $true
$r2.Font.Size = 14

# ---------------------------------------------------------------
# Edit 3: "PROGRAM " + "10" (two runs) -> "PROGRAM 10" (one run)
# ---------------------------------------------------------------
$d.Content.Find.Execute("PROGRAM 10", $true, $false, $false, $false, $false, `
    $true, 1, $false, "PROGRAM 10", 2)
